# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stocks) worksheet gets three new trailing columns:
#   H = date              (constant "2013-12-13" for every data row)
#   I = legislator_name   (constant "蘇清泉" for every data row)
#   J = legislator_id     (constant 1765 for every data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- header row -----------------------------------------------------------
$ws.Cells.Item(1,8).Value  = "date"
$ws.Cells.Item(1,9).Value  = "legislator_name"
$ws.Cells.Item(1,10).Value = "legislator_id"

# --- data rows (2 through 14) ---------------------------------------------
$firstRow = 2
$lastRow  = 14

# Assigning the literal string "2013-12-13" straight to .Value lets Excel's
# smart-typing turn it into a date serial number, which is not what the
# source data looked like (it was written as plain text). Going through a
# TEXT() formula and then collapsing the formula to its computed value with
# PasteSpecial(xlPasteValues) keeps the cell as literal text instead.
$dateRange = $ws.Range("H$firstRow`:H$lastRow")
$dateRange.Formula = "=TEXT(""2013-12-13"",""@"")"
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r,9).Value  = "蘇清泉"
    $ws.Cells.Item($r,10).Value = 1765
}
